# Update the "Statistics" sheet (sheet1): refresh existing rows 2-7 with new
# timestamps/values and append new rows 8-12.
$wsStats = $excel.ActiveWorkbook.Worksheets.Item("Statistics")

$statsData = @(
    @{ Row = 2;  Time = "2024-06-21 12:43:00"; Speed = 92.17482069321431; Density = 10 },
    @{ Row = 3;  Time = "2024-06-21 12:43:02"; Speed = 95.04422199178917; Density = 15 },
    @{ Row = 4;  Time = "2024-06-21 12:43:04"; Speed = 93.23202087367801; Density = 21 },
    @{ Row = 5;  Time = "2024-06-21 12:43:06"; Speed = 88.48835411854395; Density = 29 },
    @{ Row = 6;  Time = "2024-06-21 12:43:08"; Speed = 89.38206578331143; Density = 35 },
    @{ Row = 7;  Time = "2024-06-21 12:43:10"; Speed = 87.6110812220857;  Density = 37 },
    @{ Row = 8;  Time = "2024-06-21 12:43:12"; Speed = 87.74086381440483; Density = 37 },
    @{ Row = 9;  Time = "2024-06-21 12:43:14"; Speed = 90.20887101482903; Density = 39 },
    @{ Row = 10; Time = "2024-06-21 12:43:16"; Speed = 92.64957610426052; Density = 35 },
    @{ Row = 11; Time = "2024-06-21 12:43:18"; Speed = 91.07388210456324; Density = 36 },
    @{ Row = 12; Time = "2024-06-21 12:43:20"; Speed = 89.36651516833891; Density = 36 }
)

foreach ($entry in $statsData) {
    $r = $entry.Row
    $wsStats.Cells.Item($r, 1).Value = $entry.Time
    $wsStats.Cells.Item($r, 2).Value = $entry.Speed
    $wsStats.Cells.Item($r, 3).Value = $entry.Density
}

# Update the "Accidents" sheet (sheet2): drop the recorded accident rows,
# leaving only the header row.
$wsAccidents = $excel.ActiveWorkbook.Worksheets.Item("Accidents")
$wsAccidents.Rows("2:4").Delete()
